# Updated CEMBA overview for ALLC
#
# 1. Bump the cached "datetimeFigureOut" date field text from 1/11/21 to
#    2/2/21 on the slide master and every slide layout.
# 2. Resize/reposition the "Produce VCF" box and update its text to
#    "Produce VCF, ALLC".
# 3. Resize the connector feeding into that box to match its new position.
# 4. Nudge the enclosing group's bounding box to match (tiny +1 EMU growth
#    caused by the box's right edge moving 1 EMU further right).

$p = $ppt.ActivePresentation

function Set-DateFieldText {
    param($shapes, [string]$text)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

# --- 1. Update date placeholders on master + all layouts ---
Set-DateFieldText $p.SlideMaster.Shapes "2/2/21"
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    Set-DateFieldText $layout.Shapes "2/2/21"
}

# --- 2/3/4. Update the diagram shapes on slide 1 ---
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)   # "Group 1", the only top-level shape on the slide

# "Rounded Rectangle 14" - the "Produce VCF" box
$rect = $g.GroupItems.Item(10)
$rect.Left = 572.8696062992126
$rect.Width = 184.17133858267715
$rect.TextFrame.TextRange.Text = "Produce VCF, ALLC"

# "Elbow Connector 55" - the connector that points at the box above
$conn = $g.GroupItems.Item(22)
$conn.Width = 184.95540

# The group's own bounding box grows by a hair because the box's right
# edge now extends 1 EMU further right than before.
$g.Width = 684.8891338582678
